$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (tab title changes from 2021-09-27 to 2021-09-28)
$ws.Name = "Through 2021-09-28"

# Row 10 (August)
$ws.Range("T10").Value = 7
$ws.Range("U10").Value = 152
$ws.Range("V10").Value = 0.044

# Row 11 (September) - label and values
$ws.Range("A11").Value = "September (through 09-28)"
$ws.Range("F11").Value = 38
$ws.Range("G11").Value = 0.0952
$ws.Range("I11").Value = 67
$ws.Range("J11").Value = 0.0694
$ws.Range("L11").Value = 48
$ws.Range("M11").Value = 0.0769
$ws.Range("N11").Value = 7
$ws.Range("O11").Value = 61
$ws.Range("P11").Value = 0.1029
$ws.Range("R11").Value = 103
$ws.Range("S11").Value = 0.0374
$ws.Range("U11").Value = 166
$ws.Range("V11").Value = 0.0119

# Row 12 (Total)
$ws.Range("F12").Value = 378
$ws.Range("G12").Value = 0.1064
$ws.Range("I12").Value = 573
$ws.Range("J12").Value = 0.0803
$ws.Range("L12").Value = 481
$ws.Range("M12").Value = 0.1125
$ws.Range("N12").Value = 43
$ws.Range("O12").Value = 374
$ws.Range("P12").Value = 0.1031
$ws.Range("R12").Value = 839
$ws.Range("S12").Value = 0.0584
$ws.Range("T12").Value = 77
$ws.Range("U12").Value = 1161
$ws.Range("V12").Value = 0.0622
